# "Generate Report for Handback" - refresh the handback/handoff timestamps
# for file 482088ec-c2de-4126-81ca-50afb0d685f7.md across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
# Column G = "Latest HO Xliff Generate Date" for the
# 482088ec-c2de-4126-81ca-50afb0d685f7.md row (row 3).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-16 12:45:04"

# --- zh-cn sheet --------------------------------------------------------
# H3 = "Correspond Handoff Datetime", K3 = "Correspond Handback DateTime"
# for the same file's zh-cn row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-16 12:44:55"
$wsZhCn.Range("K3").Value = "2016-08-16 12:45:31"

# --- de-de sheet ----------------------------------------------------------
# K3 = "Correspond Handback DateTime" for the same file's de-de row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-16 12:45:38"

Write-Output "Updated handback timestamps for 482088ec-c2de-4126-81ca-50afb0d685f7.md"
